# NeapWorkflow GET_MarineTerrestrial_Extant-2020_250m.xlsx
# "updated ALUM path in workflows"
#
# The raw ABARES land-use path stored in B7 had a duplicated path segment
# (...prerelease2_20240724\ABARES_Land_use_of_Australia_2020_21_prerelease2_20240724\...).
# Fix it to the corrected single-segment path, turn it into a hyperlink (like
# the other raw/overlay paths in the sheet), and refresh the selection /
# column sizing that Excel updates as a side effect of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$correctedPath = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\inputs\raw\Land_use_of_Australia\ABARES_Land_use_of_Australia_2020_21_prerelease2_20240724\NLUM_v7p2_ALUMV8_250m_2020_21_alb.tif"
$correctedTarget = "file:///\\fs1-cbr.nexus.csiro.au\%7bev-neap%7d\work\extent\inputs\raw\Land_use_of_Australia\ABARES_Land_use_of_Australia_2020_21_prerelease2_20240724\NLUM_v7p2_ALUMV8_250m_2020_21_alb.tif"

$cellB7 = $ws.Range("B7")
$cellB7.Value = $correctedPath

$ws.Hyperlinks.Add($cellB7, $correctedTarget)
$cellB7.Style = "Hyperlink"

# Column B now needs to be wider to fit the (still long) corrected path.
$ws.Columns.Item(2).ColumnWidth = 220.5

# Reflect the author's final cursor position in the sheet.
[void]$ws.Range("B15").Select()
